$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date (serial 43811 -> 2019-12-12)
$ws.Range("B3").Value = [DateTime]"2019-12-12"

# Team name / numeric entries
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 5

# Team member salaries
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 100

# Update selection to match target workbook
$ws.Range("B6").Select()
